# Generate Report for Handoff
# Replaces the two tracked e2e files with a freshly-handed-off pair:
#   0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md  -> 0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md
#   e3817207-2a06-43a7-a18f-c31df4325e01.md  -> ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md
# and flips their status from "handed back" to "ready for handoff", updating
# the handoff file / datetime columns accordingly (both localized sheets
# now reference a single freshly generated xliff, with row 2 = primary and
# row 3 flagged as a content duplicate of it).

$wb = $excel.ActiveWorkbook

$oldFile1 = "0a34ea6d-c30a-4cf0-a4b2-c4e7a97d0f6b.md"
$newFile1 = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.md"
$oldFile2 = "e3817207-2a06-43a7-a18f-c31df4325e01.md"
$newFile2 = "ffff2fd1e58b-2408-4bc1-9197-ecfdddf93ad2.md"

$newXlfZh = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.zh-cn.xlf"
$newXlfDe = "0ffdf4c3-554d-4d63-a026-6e8bfa94ab7d.e24154a4794c8801e1eece2bf03a3390907cfa15.de-de.xlf"

$statusText = "Ready for handoff"
$handoffDateZh = "2016-08-30 01:04:35"
$handoffDateDe = "2016-08-30 01:04:40"
$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = "e2e\$newFile1"
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("G2").Value = $handoffDateDe

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = "e2e\$newFile2"
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $handoffDateDe

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newFile1"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newFile2"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $handoffDateZh
$wsZh.Range("K2").Value = $zeroDate

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $handoffDateZh
$wsZh.Range("K3").Value = $zeroDate

# I2/I3 (Latest Target File) and J2/J3 (Latest Handback File) lose their
# values - no handback has happened yet for the freshly generated handoff.
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I2").Value = "'"
$wsZh.Range("J2").Value = "'"

$wsZh.Range("I3").Style = "Normal"
$wsZh.Range("I3").Value = "'"
$wsZh.Range("J3").Value = "'"

$toDeleteZh = @()
foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    } elseif (($addr -eq '$I$2') -or ($addr -eq '$I$3')) {
        $toDeleteZh += $hl
    }
}
foreach ($hl in $toDeleteZh) {
    $hl.Delete()
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $handoffDateDe
$wsDe.Range("K2").Value = $zeroDate

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $handoffDateDe
$wsDe.Range("K3").Value = $zeroDate

$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I2").Value = "'"
$wsDe.Range("J2").Value = "'"

$wsDe.Range("I3").Style = "Normal"
$wsDe.Range("I3").Value = "'"
$wsDe.Range("J3").Value = "'"

$toDeleteDe = @()
foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    } elseif (($addr -eq '$I$2') -or ($addr -eq '$I$3')) {
        $toDeleteDe += $hl
    }
}
foreach ($hl in $toDeleteDe) {
    $hl.Delete()
}
